# progressbypushkar.xlsx - "Add files via upload"
# Adds a new progress-log entry (row 5: date + work done) to Sheet1,
# widens the DATE/WORK/(spare) columns to fit the new content, and
# moves the active selection to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row -----------------------------------------------------
# A5: the date of the new entry (centered, like the other DATE cells)
$ws.Range("A5").Value = "31-1-2019"
$ws.Range("A5").HorizontalAlignment = -4108   # xlCenter

# B5: the work done that day (plain/general style)
$ws.Range("B5").Value = "Changed pom.xml,Created Student Repo,Student Model,Controller Error running TestApplication.java"

# --- Column widths ------------------------------------------------------
# Column A (DATE) stays narrow, column B (WORK) grows very wide to fit
# the long description, column C gets a bit of extra room too.
$ws.Columns.Item(1).ColumnWidth = 10.5
$ws.Columns.Item(2).ColumnWidth = 81.5
$ws.Columns.Item(3).ColumnWidth = 24.166666666666668

# --- Selection ------------------------------------------------------
$ws.Range("B4").Select() | Out-Null
